# Simulated Wild Card round and logged it
# Appends the play-by-play logs for the extra playoff game to the YDS and
# ST per-play strings, and bumps the aggregate stat totals on the
# OFF / DEF / ST / TURNS / PEN summary sheets to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# YDS sheet: append the new game's rush/pass yardage logs
# ---------------------------------------------------------------
$wsYDS = $wb.Worksheets.Item("YDS")

$wsYDS.Range("B2").Value = $wsYDS.Range("B2").Value() + " 3 2 6 2 6 5 4 2 -2 6 2 1 0 1 9 0 8 2 19 1 8 11 5 -4 9 10 2 -1 2 3 11 10 0 2 2 4 3 2 2 4 0 3 5"
$wsYDS.Range("B3").Value = $wsYDS.Range("B3").Value() + " 1 0 13 5 7 4 2 4 5 8 1 8 13 9 8 35 14 20 19 10 6 9 9 15 -4 9 10 0 2 -4 6 5 11 -2 6 11 0 7 23 5 7 4 8 13 12 13 6"
$wsYDS.Range("C2").Value = $wsYDS.Range("C2").Value() + " 4 20 6 10 0 2 -2 2 2 0 0 7 8 4 4 3 3 2 1 0 7 29 5 1 4 2 62 9 8 35 0 7 2 1 17 6 1 1 -2 15 5 3 1 0 1 4 2 -2 8 8 -3 3 4 4 8 11 10 -4 1 7 2 3 -2 6 3 3 -3"
$wsYDS.Range("C3").Value = $wsYDS.Range("C3").Value() + " 2 43 9 7 23 31 12 5 22 17 20 9 2 25 11 5 17 22 4 14 31 10 2 40 4 7 41 29 16"

# ---------------------------------------------------------------
# OFF sheet: bump offensive aggregate totals
# ---------------------------------------------------------------
$wsOFF = $wb.Worksheets.Item("OFF")

$wsOFF.Range("B2").Value = 6
$wsOFF.Range("C2").Value = 240
$wsOFF.Range("D2").Value = 14
$wsOFF.Range("E2").Value = 12
$wsOFF.Range("F2").Value = 44
$wsOFF.Range("G2").Value = 73
$wsOFF.Range("H2").Value = 10
$wsOFF.Range("I2").Value = 15
$wsOFF.Range("J2").Value = 40
$wsOFF.Range("L2").Value = 287
$wsOFF.Range("M2").Value = 196
$wsOFF.Range("N2").Value = 17
$wsOFF.Range("O2").Value = 29
$wsOFF.Range("P2").Value = 17
$wsOFF.Range("Q2").Value = 517

$wsOFF.Range("B3").Value = 13
$wsOFF.Range("C3").Value = 193
$wsOFF.Range("E3").Value = 45
$wsOFF.Range("F3").Value = 145
$wsOFF.Range("G3").Value = 52
$wsOFF.Range("H3").Value = 36
$wsOFF.Range("I3").Value = 55
$wsOFF.Range("J3").Value = 60
$wsOFF.Range("L3").Value = 329
$wsOFF.Range("M3").Value = 231
$wsOFF.Range("N3").Value = 24
$wsOFF.Range("Q3").Value = 665

# ---------------------------------------------------------------
# DEF sheet: bump defensive aggregate totals
# ---------------------------------------------------------------
$wsDEF = $wb.Worksheets.Item("DEF")

$wsDEF.Range("B2").Value = 8
$wsDEF.Range("C2").Value = 221
$wsDEF.Range("D2").Value = 13
$wsDEF.Range("E2").Value = 16
$wsDEF.Range("F2").Value = 72
$wsDEF.Range("G2").Value = 57
$wsDEF.Range("H2").Value = 6
$wsDEF.Range("I2").Value = 7
$wsDEF.Range("J2").Value = 33
$wsDEF.Range("L2").Value = 248
$wsDEF.Range("M2").Value = 165
$wsDEF.Range("N2").Value = 30
$wsDEF.Range("Q2").Value = 495

$wsDEF.Range("C3").Value = 181
$wsDEF.Range("D3").Value = 6
$wsDEF.Range("E3").Value = 39
$wsDEF.Range("F3").Value = 104
$wsDEF.Range("G3").Value = 36
$wsDEF.Range("H3").Value = 32
$wsDEF.Range("I3").Value = 57
$wsDEF.Range("J3").Value = 65
$wsDEF.Range("L3").Value = 332
$wsDEF.Range("M3").Value = 217
$wsDEF.Range("N3").Value = 13
$wsDEF.Range("Q3").Value = 623

# ---------------------------------------------------------------
# ST sheet: append kickoff/punt distance & return logs, bump totals
# ---------------------------------------------------------------
$wsST = $wb.Worksheets.Item("ST")

$wsST.Range("B4").Value = $wsST.Range("B4").Value() + " 70 61"
$wsST.Range("B5").Value = $wsST.Range("B5").Value() + " 21 7"
$wsST.Range("B6").Value = $wsST.Range("B6").Value() + " 6 14 15"
$wsST.Range("D3").Value = $wsST.Range("D3").Value() + " 57 49 48 51 52 54 37 47"
$wsST.Range("D4").Value = $wsST.Range("D4").Value() + " 13 0 7 0 8 0 11 0"
$wsST.Range("D5").Value = $wsST.Range("D5").Value() + " 0 0 0 0 0 0"

$wsST.Range("B2").Value = 96
$wsST.Range("D2").Value = 57
$wsST.Range("F2").Value = 552
$wsST.Range("G2").Value = 538
$wsST.Range("H2").Value = 4
$wsST.Range("I2").Value = 1
$wsST.Range("J2").Value = 217
$wsST.Range("K2").Value = 204
$wsST.Range("L2").Value = 117
$wsST.Range("M2").Value = 81
$wsST.Range("N2").Value = 87
$wsST.Range("O2").Value = 66

$wsST.Range("B3").Value = 74

# ---------------------------------------------------------------
# TURNS sheet: bump turnover totals
# ---------------------------------------------------------------
$wsTURNS = $wb.Worksheets.Item("TURNS")

$wsTURNS.Range("C2").Value = 3
$wsTURNS.Range("E2").Value = 8

$wsTURNS.Range("B3").Value = 5
$wsTURNS.Range("E3").Value = 8

# ---------------------------------------------------------------
# PEN sheet: bump penalty totals
# ---------------------------------------------------------------
$wsPEN = $wb.Worksheets.Item("PEN")

$wsPEN.Range("B2").Value = 28
$wsPEN.Range("D2").Value = 8

$wsPEN.Range("B3").Value = 22
$wsPEN.Range("D3").Value = 5

$wsPEN.Range("B4").Value = 4
